$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 842.0833
$ws.Range("J32").Value = 905.4
$ws.Range("L32").Value = 905.4
$ws.Range("N32").Value = -1557.4
$ws.Range("H33").Value = 244
$ws.Range("I33").Value = 206
$ws.Range("J33").Value = 434
$ws.Range("K33").Value = 206
$ws.Range("L33").Value = 434
$ws.Range("M33").Value = 23
$ws.Range("N33").Value = -892
$ws.Range("H98").Value = 820.1177
$ws.Range("I98").Value = 809
$ws.Range("J98").Value = 998
$ws.Range("K98").Value = 809
$ws.Range("L98").Value = 998
$ws.Range("M98").Value = 689
$ws.Range("N98").Value = -3994
$ws.Range("H122").Value = 820.1177
$ws.Range("I122").Value = 809
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 2427
$ws.Range("L122").Value = 2994
$ws.Range("M122").Value = 23
$ws.Range("N122").Value = -7894

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1168.0454
$ws.Range("I45").Value = 1081.0625
$ws.Range("J45").Value = 1400
$ws.Range("K45").Value = 1081.0625
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = -704.0625
$ws.Range("N45").Value = -2154
$ws.Range("H61").Value = 125252130
$ws.Range("I61").Value = 250252880
$ws.Range("J61").Value = 251373.5
$ws.Range("K61").Value = 250252880
$ws.Range("L61").Value = 251373.5
$ws.Range("M61").Value = -250252668
$ws.Range("N61").Value = -251797.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H132").Value = 9657500
$ws.Range("I132").Value = 12847767
$ws.Range("J132").Value = 86698.30499999999
$ws.Range("K132").Value = 38543301
$ws.Range("L132").Value = 260094.915
$ws.Range("M132").Value = -38540771
$ws.Range("N132").Value = -265154.915
$ws.Range("H136").Value = 125252130
$ws.Range("I136").Value = 250252880
$ws.Range("J136").Value = 251373.5
$ws.Range("K136").Value = 750758640
$ws.Range("L136").Value = 754120.5
$ws.Range("M136").Value = -750756090
$ws.Range("N136").Value = -759220.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 354.66666
$ws.Range("I94").Value = 372.125
$ws.Range("J94").Value = 334.7143
$ws.Range("K94").Value = 372.125
$ws.Range("L94").Value = 334.7143
$ws.Range("M94").Value = 78.875
$ws.Range("N94").Value = -1236.7143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 25082.094
$ws.Range("I132").Value = 1532.25
$ws.Range("J132").Value = 93590.73
$ws.Range("K132").Value = 4596.75
$ws.Range("L132").Value = 280772.19
$ws.Range("M132").Value = -2066.75
$ws.Range("N132").Value = -285832.19

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 600.0625
$ws.Range("I5").Value = 426.18518
$ws.Range("J5").Value = 1539
$ws.Range("K5").Value = 1278.55554
$ws.Range("L5").Value = 4617
$ws.Range("M5").Value = -1166.55554
$ws.Range("N5").Value = -4841
$ws.Range("H92").Value = 900.4666999999999
$ws.Range("I92").Value = 900.4666999999999
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2701.4001
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1453.4001
$ws.Range("N92").ClearContents()
$ws.Range("H122").Value = 828.0833
$ws.Range("I122").Value = 198.85715
$ws.Range("J122").Value = 1228.5
$ws.Range("K122").Value = 1789.71435
$ws.Range("L122").Value = 11056.5
$ws.Range("M122").Value = 660.28565
$ws.Range("N122").Value = -15956.5
$ws.Range("H123").Value = 1889.2858
$ws.Range("I123").Value = 832.5
$ws.Range("J123").Value = 3298.3333
$ws.Range("K123").Value = 2497.5
$ws.Range("L123").Value = 9894.999899999999
$ws.Range("M123").Value = -47.5
$ws.Range("N123").Value = -14794.9999
$ws.Range("H125").Value = 2970
$ws.Range("I125").Value = 1910
$ws.Range("J125").Value = 3500
$ws.Range("K125").Value = 5730
$ws.Range("L125").Value = 10500
$ws.Range("M125").Value = -810
$ws.Range("N125").Value = -20340
$ws.Range("H131").Value = 908.0328
$ws.Range("I131").Value = 499.5
$ws.Range("J131").Value = 936.7018
$ws.Range("K131").Value = 1498.5
$ws.Range("L131").Value = 2810.1054
$ws.Range("M131").Value = 3541.5
$ws.Range("N131").Value = -12890.1054
$ws.Range("H132").Value = 3435.842
$ws.Range("J132").Value = 4846.5
$ws.Range("L132").Value = 43618.5
$ws.Range("N132").Value = -48678.5
$ws.Range("H133").Value = 3000
$ws.Range("I133").Value = 3000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 9000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -3940
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 6851.8423
$ws.Range("I134").Value = 4912.143
$ws.Range("J134").Value = 7983.3335
$ws.Range("K134").Value = 14736.429
$ws.Range("L134").Value = 23950.0005
$ws.Range("M134").Value = -9666.429
$ws.Range("N134").Value = -34090.00049999999
$ws.Range("H135").Value = 600.0625
$ws.Range("I135").Value = 426.18518
$ws.Range("J135").Value = 1539
$ws.Range("K135").Value = 3835.66662
$ws.Range("L135").Value = 13851
$ws.Range("M135").Value = -1300.66662
$ws.Range("N135").Value = -18921
$ws.Range("H136").Value = 3260.0527
$ws.Range("I136").Value = 2878.6667
$ws.Range("J136").Value = 3603.3
$ws.Range("K136").Value = 8636.000100000001
$ws.Range("L136").Value = 10809.9
$ws.Range("M136").Value = -3536.000100000001
$ws.Range("N136").Value = -21009.9
$ws.Range("H137").Value = 1857
$ws.Range("I137").Value = 956
$ws.Range("J137").Value = 2758
$ws.Range("K137").Value = 2868
$ws.Range("L137").Value = 8274
$ws.Range("M137").Value = 2232
$ws.Range("N137").Value = -18474
$ws.Range("H139").Value = 4394.1
$ws.Range("I139").Value = 2000.4375
$ws.Range("J139").Value = 7129.7144
$ws.Range("K139").Value = 6001.3125
$ws.Range("L139").Value = 21389.1432
$ws.Range("M139").Value = -861.3125
$ws.Range("N139").Value = -31669.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 33465.57
$ws.Range("I70").Value = 48178.13
$ws.Range("J70").Value = 5266.5
$ws.Range("K70").Value = 48178.13
$ws.Range("L70").Value = 5266.5
$ws.Range("M70").Value = -47908.13
$ws.Range("N70").Value = -5806.5
$ws.Range("H73").Value = 33465.57
$ws.Range("I73").Value = 48178.13
$ws.Range("J73").Value = 5266.5
$ws.Range("K73").Value = 48178.13
$ws.Range("L73").Value = 5266.5
$ws.Range("M73").Value = -47242.13
$ws.Range("N73").Value = -7138.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 40148.8
$ws.Range("J133").Value = 40148.8
$ws.Range("L133").Value = 40148.8
$ws.Range("N133").Value = -45208.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1347.0834
$ws.Range("I126").Value = 683.125
$ws.Range("J126").Value = 2675
$ws.Range("K126").Value = 2049.375
$ws.Range("L126").Value = 8025
$ws.Range("M126").Value = 420.625
$ws.Range("N126").Value = -12965
